# Auto-generated Excel COM-interop script to apply crypto price/volume updates
# per commit "Updated cryptos list on Sun Feb 11 12:51:25 UTC 2024 with GitHub Actions"
#
# Note: several Price (D) values are plain decimals (e.g. "1.00", "9.12") that
# Excel's COM layer would otherwise auto-coerce into numeric cells. The source
# workbook stores every Price/Volume cell as text, so those assignments are
# prefixed with a leading apostrophe (Excel's standard 'force text' entry) to
# keep them text after the write, matching the original cell typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "48.338.44"
$ws.Range("E2").Value = "  +2.41%  "

# Row 3
$ws.Range("D3").Value = "2.524.84"
$ws.Range("E3").Value = "  +1.58%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'110.03"
$ws.Range("E5").Value = "  +2.31%  "

# Row 6
$ws.Range("D6").Value = "'322.33"
$ws.Range("E6").Value = "  +0.45%  "

# Row 7
$ws.Range("D7").Value = "'0.532"
$ws.Range("E7").Value = "  +2.32%  "

# Row 8
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("E9").Value = "  +3.72%  "

# Row 10
$ws.Range("D10").Value = "'40.63"
$ws.Range("E10").Value = "  +5.74%  "

# Row 11
$ws.Range("D11").Value = "'20.45"
$ws.Range("E11").Value = "  +12.32%  "

# Row 12
$ws.Range("E12").Value = "  +1.82%  "

# Row 13
$ws.Range("E13").Value = "  +1.20%  "

# Row 14
$ws.Range("D14").Value = "'7.28"
$ws.Range("E14").Value = "  +2.64%  "

# Row 15
$ws.Range("D15").Value = "2.922.04"
$ws.Range("E15").Value = "  +1.70%  "

# Row 16
$ws.Range("D16").Value = "2.525.26"
$ws.Range("E16").Value = "  +1.51%  "

# Row 17
$ws.Range("D17").Value = "'0.854"
$ws.Range("E17").Value = "  +1.19%  "

# Row 18
$ws.Range("D18").Value = "48.172.22"
$ws.Range("E18").Value = "  +2.24%  "

# Row 19
$ws.Range("D19").Value = "'13.35"
$ws.Range("E19").Value = "  +5.02%  "

# Row 20
$ws.Range("D20").Value = "'6.63"
$ws.Range("E20").Value = "  +0.33%  "

# Row 21
$ws.Range("E21").Value = "  +2.11%  "

# Row 22
$ws.Range("D22").Value = "'2.68"
$ws.Range("E22").Value = "  -0.94%  "

# Row 23
$ws.Range("D23").Value = "'71.99"
$ws.Range("E23").Value = "  +2.54%  "

# Row 24
$ws.Range("D24").Value = "'270.47"
$ws.Range("E24").Value = "  +10.27%  "

# Row 25
$ws.Range("E25").Value = "  +0.92%  "

# Row 26
$ws.Range("E26").Value = "  +0.01%  "

# Row 27
$ws.Range("D27").Value = "'26.04"
$ws.Range("E27").Value = "  +1.60%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.28"
$ws.Range("E28").Value = "  +0.13%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'10.14"
$ws.Range("E29").Value = "  +1.43%  "

# Row 30
$ws.Range("E30").Value = "  +7.31%  "

# Row 31
$ws.Range("D31").Value = "'35.65"
$ws.Range("E31").Value = "  +3.90%  "

# Row 32
$ws.Range("D32").Value = "'49.67"
$ws.Range("E32").Value = "  +0.41%  "

# Row 33
$ws.Range("D33").Value = "'19.72"
$ws.Range("E33").Value = "  -2.51%  "

# Row 34
$ws.Range("D34").Value = "'5.39"
$ws.Range("E34").Value = "  +1.41%  "

# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.01%  "

# Row 36
$ws.Range("E36").Value = "  +1.23%  "

# Row 37
$ws.Range("E37").Value = "  +1.64%  "

# Row 38
$ws.Range("D38").Value = "'4.69"
$ws.Range("E38").Value = "  +2.11%  "

# Row 39
$ws.Range("D39").Value = "'3.02"

# Row 40
$ws.Range("E40").Value = "  +0.79%  "

# Row 41
$ws.Range("D41").Value = "'121.85"
$ws.Range("E41").Value = "  +2.57%  "

# Row 42
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'2.21"
$ws.Range("E42").Value = "  -0.27%  "

# Row 43
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'21.87"
$ws.Range("E43").Value = "  -1.66%  "

# Row 44
$ws.Range("E44").Value = "  +2.49%  "

# Row 45
$ws.Range("D45").Value = "2.027.01"
$ws.Range("E45").Value = "  +2.24%  "

# Row 46
$ws.Range("D46").Value = "'3.15"
$ws.Range("E46").Value = "  +5.32%  "

# Row 47
$ws.Range("D47").Value = "'1.89"
$ws.Range("E47").Value = "  +7.91%  "

# Row 48
$ws.Range("E48").Value = "  +2.69%  "

# Row 49
$ws.Range("D49").Value = "'9.12"
$ws.Range("E49").Value = "  +1.17%  "

# Row 50
$ws.Range("D50").Value = "'5.22"
$ws.Range("E50").Value = "  +2.29%  "

# Row 51
$ws.Range("D51").Value = "'79.53"
$ws.Range("E51").Value = "  +3.66%  "

